# Updates the cryptocurrency price/volume table to the latest scraped values.
# (Mirrors the GitHub Actions "update cryptos list" job output.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "96.455.62"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.707.66"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.16"
$ws.Range("E5").Value = "  -3.55%  "
$ws.Range("E6").Value = "  -2.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "651.04"
$ws.Range("E7").Value = "  -3.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.427"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  -7.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.707.60"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("B12").Value = "ShibaInu"
$ws.Range("C12").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000309"
$ws.Range("E12").Value = "  +14.83%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "44.28"
$ws.Range("E13").Value = "  -2.62%  "
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.397.68"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.153.90"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.83"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.716.43"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.14"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.503"
$ws.Range("E22").Value = "  -8.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "520.24"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.39"
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.45"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.32"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("E29").Value = "  +6.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.00"
$ws.Range("E30").Value = "  -3.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.14"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.87"
$ws.Range("E33").Value = "  +7.05%  "
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.27"
$ws.Range("E36").Value = "  -3.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "644.75"
$ws.Range("E37").Value = "  +4.37%  "
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("E39").Value = "  -0.63%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.82"
$ws.Range("E41").Value = "  +10.66%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "41.04"
$ws.Range("E42").Value = "  -3.62%  "
$ws.Range("E43").Value = "  +3.29%  "
$ws.Range("E44").Value = "  -2.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.966"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0454"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.434"
$ws.Range("E47").Value = "  +2.15%  "
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.48"
$ws.Range("E50").Value = "  -1.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.50"
$ws.Range("E51").Value = "  +1.05%  "
